# Auto-generated edit script applying the Adamantoise_Profits market-data refresh
# (currentAveragePrice / LevePrice / LeveProfit columns) across the 7 affected sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 11
$ws.Range("H11").Value = 27778120
$ws.Range("I11").Value = 27778120
$ws.Range("K11").Value = 27778120
$ws.Range("M11").Value = -27777980
# Row 53
$ws.Range("H53").Value = 1002.2692
$ws.Range("I53").Value = 676.38464
$ws.Range("J53").Value = 1328.1538
$ws.Range("K53").Value = 676.38464
$ws.Range("L53").Value = 1328.1538
$ws.Range("M53").Value = -39.38463999999999
$ws.Range("N53").Value = -2602.1538
# Row 88
$ws.Range("H88").Value = 1359.3636
$ws.Range("I88").Value = 1673
$ws.Range("J88").Value = 1045.7273
$ws.Range("K88").Value = 1673
$ws.Range("L88").Value = 1045.7273
$ws.Range("M88").Value = -1267
$ws.Range("N88").Value = -1857.7273
# Row 91
$ws.Range("H91").Value = 1359.3636
$ws.Range("I91").Value = 1673
$ws.Range("J91").Value = 1045.7273
$ws.Range("K91").Value = 1673
$ws.Range("L91").Value = 1045.7273
$ws.Range("M91").Value = -269
$ws.Range("N91").Value = -3853.7273
# Row 93
$ws.Range("H93").Value = 69989.664
$ws.Range("J93").Value = 69989.664
$ws.Range("L93").Value = 69989.664
$ws.Range("N93").Value = -74981.664
# Row 116
$ws.Range("H116").Value = 11224.85
$ws.Range("I116").Value = 11976.883
$ws.Range("J116").Value = 6963.3335
$ws.Range("K116").Value = 11976.883
$ws.Range("L116").Value = 6963.3335
$ws.Range("M116").Value = -8534.883
$ws.Range("N116").Value = -13847.3335
# Row 118
$ws.Range("H118").Value = 1672.7
$ws.Range("I118").Value = 995.5714
$ws.Range("K118").Value = 2986.7142
$ws.Range("M118").Value = -1329.7142
# Row 132
$ws.Range("H132").Value = 5610.8623
$ws.Range("I132").Value = 6409.8184
$ws.Range("K132").Value = 19229.4552
$ws.Range("M132").Value = -16699.4552
# Row 135
$ws.Range("H135").Value = 1302.5
$ws.Range("I135").Value = 1302.5
$ws.Range("K135").Value = 11722.5
$ws.Range("M135").Value = -9187.5

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 4471.273
$ws.Range("I32").Value = 1256.9706
$ws.Range("J32").Value = 15399.9
$ws.Range("K32").Value = 1256.9706
$ws.Range("L32").Value = 15399.9
$ws.Range("M32").Value = -969.9706000000001
$ws.Range("N32").Value = -15973.9
# Row 61
$ws.Range("H61").Value = 2931.4285
$ws.Range("I61").Value = 2783
$ws.Range("K61").Value = 2783
$ws.Range("M61").Value = -2571
# Row 63
$ws.Range("H63").Value = 3829.8
$ws.Range("I63").Value = 2299.6
$ws.Range("J63").Value = 5360
$ws.Range("K63").Value = 2299.6
$ws.Range("L63").Value = 5360
$ws.Range("M63").Value = -1613.6
$ws.Range("N63").Value = -6732
# Row 66
$ws.Range("H66").Value = 3829.8
$ws.Range("I66").Value = 2299.6
$ws.Range("J66").Value = 5360
$ws.Range("K66").Value = 11498
$ws.Range("L66").Value = 26800
$ws.Range("M66").Value = -8066
$ws.Range("N66").Value = -33664
# Row 103
$ws.Range("H103").Value = 68804.5
$ws.Range("J103").Value = 68804.5
$ws.Range("L103").Value = 68804.5
$ws.Range("N103").Value = -71148.5
# Row 122
$ws.Range("H122").Value = 6229.311
$ws.Range("I122").Value = 4582.3213
$ws.Range("K122").Value = 13746.9639
$ws.Range("M122").Value = -11296.9639
# Row 132
$ws.Range("H132").Value = 335878.38
$ws.Range("I132").Value = 347287.47
$ws.Range("K132").Value = 1041862.41
$ws.Range("M132").Value = -1039332.41
# Row 136
$ws.Range("H136").Value = 2931.4285
$ws.Range("I136").Value = 2783
$ws.Range("K136").Value = 8349
$ws.Range("M136").Value = -5799

$ws = $wb.Worksheets.Item("BSM")
# Row 56
$ws.Range("H56").Value = 22855.715
# Row 134
$ws.Range("H134").Value = 9805881
$ws.Range("I134").Value = 1607.3793
$ws.Range("K134").Value = 4822.1379
$ws.Range("M134").Value = -2287.1379

$ws = $wb.Worksheets.Item("CUL")
# Row 40
$ws.Range("H40").Value = 1036.1904
$ws.Range("J40").Value = 5037.5
$ws.Range("L40").Value = 20150
$ws.Range("N40").Value = -20288
# Row 131
$ws.Range("H131").Value = 1495.2565
$ws.Range("J131").Value = 1958.8948
$ws.Range("L131").Value = 5876.6844
$ws.Range("N131").Value = -15956.6844

$ws = $wb.Worksheets.Item("GSM")
# Row 43
$ws.Range("H43").Value = 280
$ws.Range("I43").Value = 280
$ws.Range("K43").Value = 280
$ws.Range("M43").Value = -129
# Row 102
$ws.Range("H102").Value = 2822.5
$ws.Range("I102").Value = 2774.2727
$ws.Range("K102").Value = 2774.2727
$ws.Range("M102").Value = -1152.2727
# Row 122
$ws.Range("H122").Value = 3875.8125
$ws.Range("I122").Value = 3626.625
$ws.Range("K122").Value = 10879.875
$ws.Range("M122").Value = -8429.875
# Row 126
$ws.Range("H126").Value = 4250
$ws.Range("J126").Value = 4457.5
$ws.Range("L126").Value = 13372.5
$ws.Range("N126").Value = -18312.5
# Row 132
$ws.Range("H132").Value = 2672.76
$ws.Range("I132").Value = 2719.1365
$ws.Range("K132").Value = 8157.4095
$ws.Range("M132").Value = -5627.4095

$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 1680.2
$ws.Range("I22").Value = 1600.25
$ws.Range("K22").Value = 1600.25
$ws.Range("M22").Value = -1305.25
# Row 27
$ws.Range("H27").Value = 1680.2
$ws.Range("I27").Value = 1600.25
$ws.Range("K27").Value = 1600.25
$ws.Range("M27").Value = -1493.25
# Row 56
$ws.Range("H56").Value = 33257.75
$ws.Range("I56").Value = 28010.334
$ws.Range("K56").Value = 28010.334
$ws.Range("M56").Value = -27319.334
# Row 132
$ws.Range("H132").Value = 591999.5
$ws.Range("I132").Value = 772353.25
$ws.Range("J132").Value = 5850
$ws.Range("K132").Value = 2317059.75
$ws.Range("L132").Value = 17550
$ws.Range("M132").Value = -2314529.75
$ws.Range("N132").Value = -22610
# Row 136
$ws.Range("H136").Value = 5150.2
$ws.Range("I136").Value = 4252
$ws.Range("K136").Value = 12756
$ws.Range("M136").Value = -10206

$ws = $wb.Worksheets.Item("WVR")
# Row 32
$ws.Range("H32").Value = 9356.5
$ws.Range("I32").Value = 9356.5
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 9356.5
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -9039.5
$ws.Range("N32").ClearContents()
# Row 62
$ws.Range("H62").Value = 4893.125
$ws.Range("J62").Value = 6375
$ws.Range("L62").Value = 6375
$ws.Range("N62").Value = -7623
# Row 65
$ws.Range("H65").Value = 4893.125
$ws.Range("J65").Value = 6375
$ws.Range("L65").Value = 31875
$ws.Range("N65").Value = -38115
# Row 81
$ws.Range("H81").Value = 103614.5
$ws.Range("I81").Value = 127018.125
$ws.Range("K81").Value = 254036.25
$ws.Range("M81").Value = -252975.25
# Row 84
$ws.Range("H84").Value = 103614.5
$ws.Range("I84").Value = 127018.125
$ws.Range("K84").Value = 1270181.25
$ws.Range("M84").Value = -1264877.25
# Row 132
$ws.Range("H132").Value = 45163.08
$ws.Range("I132").Value = 59805.89
$ws.Range("J132").Value = 7510.143
$ws.Range("K132").Value = 179417.67
$ws.Range("L132").Value = 22530.429
$ws.Range("M132").Value = -176887.67
$ws.Range("N132").Value = -27590.429
# Row 136
$ws.Range("H136").Value = 25783.512
$ws.Range("I136").Value = 1859
$ws.Range("K136").Value = 5577
$ws.Range("M136").Value = -3027
